$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Russian geometry names in column A (data rows 2-417) with the
# English equivalents ("в merge with spaceclaim" commit: localize geometry
# labels). Three contiguous blocks, one per geometry type.
$ws.Range("A2:A113").Value()   = "cube"
$ws.Range("A114:A225").Value() = "sphere"
$ws.Range("A226:A417").Value() = "ellipse"

# Scroll the frozen view down and move the selection to where the user was
# last working.
$win = $excel.ActiveWindow
$win.ScrollRow = 207
$win.ScrollColumn = 1
$ws.Range("J218").Select()
